$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Lightning Ninjutsu Scroll[Ninjutsu Scroll]"
$ws.Range("B2").Value = "Stability %`n5Katana only:Accuracy %`n10"

$ws.Range("A3").Value = "Water Ninjutsu Scroll[Ninjutsu Scroll]"
$ws.Range("B3").Value = "Ailment Resistance %`n5Magic Device only:`nAggro %`n-10"

$ws.Range("A4").Value = "Earth Ninjutsu Scroll[Ninjutsu Scroll]"
$ws.Range("B4").Value = "`nMaxHP %`n101-Handed Sword only:Fractional `nBarrier %`n10"

$ws.Range("A5").Value = "Metal Ninjutsu Scroll[Ninjutsu Scroll]"
$ws.Range("B5").Value = "`nCritical Rate`n5"

$ws.Range("A6").Value = "Fire Ninjutsu Scroll[Ninjutsu Scroll]"
$ws.Range("B6").Value = "MATK`n %`n1Staff only:Magic Pierce %`n5"

$ws.Range("A7").Value = "Wind Ninjutsu Scroll[Ninjutsu Scroll]"
$ws.Range("B7").Value = "`nASPD250Katana only:`nCritical Rate`n5"

$ws.Range("A8").Value = "Dark Ninjutsu Scroll[Ninjutsu Scroll]"
$ws.Range("B8").Value = "`nAggro %`n-10"
